# Enhance Setup table decoupling and improve financial linkage logic
$wb = $excel.ActiveWorkbook

# --- Setup sheet: decouple the "Credit Card" row's Sub-Category from its
#     Type so it reflects the CC account config key instead of duplicating
#     the Liability type. ---
$setup = $wb.Worksheets.Item("Setup")
$setup.Range("B5").Value = "CC"

# --- VERSION sheet: record the re-generated build/version metadata and a
#     fresh import-history entry for the latest (re-)run. ---
$version = $wb.Worksheets.Item("VERSION")
$version.Range("B2").Value = "1.1.0"
$version.Range("B3").Value = "d30caf8818b4bcec55d88484ee588c128f9dcf14"
$version.Range("B4").Value = "12/30/2025, 11:11:02 AM"
$version.Range("A7").Value = "Import at 12/30/2025, 11:11:02 AM"
$version.Range("A8").Value = "Import at 12/30/2025, 11:11:03 AM"
